# v1.2.0: add a new "ton420ls" parameter sheet (copy of "ton50ls") with its
# own set of values, placed before the original "ton50ls" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Locate the existing "ton50ls" sheet -------------------------------
$ton50 = $wb.Worksheets.Item("ton50ls")

# --- 2. Duplicate it, inserting the copy *before* ton50ls -----------------
#     (sheet references here are position-based, so re-resolve by name
#      afterwards instead of trusting stale variables)
$ton50.Copy($ton50)

$ton420 = $wb.Worksheets.Item(1)
$ton420.Name = "ton420ls"
$ton50 = $wb.Worksheets.Item("ton50ls")

# --- 3. Recreate the data table lost by the copy on the new sheet ---------
$tbl = $ton420.ListObjects.Add(1, $ton420.Range("A1:C34"), $null, 1)
$tbl.Name = "Tabla13"
$tbl.TableStyle = "TableStyleMedium11"

# --- 4. Update the C-column values for the ton420ls parameter set ---------
$values = @{
    2  = 6000
    3  = 12000
    4  = 18000
    5  = 21500
    6  = 1000
    7  = 500
    8  = 100
    9  = 50
    10 = 21500
    11 = 18000
    12 = 12000
    13 = 6000
    14 = 20
    15 = 150
    16 = 750
    17 = 1500
    18 = 1
    19 = 0
    20 = 50
    21 = 225
    22 = 500
    23 = 20
    24 = 400
    25 = 10
    26 = 10
    27 = 1
    28 = 2
    29 = 100
    30 = 5
    31 = 10
    32 = 5
    33 = 1
    34 = 60
}

foreach ($row in $values.Keys) {
    $ton420.Cells.Item($row, 3).Value = $values[$row]
}

# --- 5. Selections & active sheet, matching the authored workbook view ----
#     (Range.Select applies to whichever sheet is active, so activate
#      each sheet before selecting on it)
$ton50.Activate()
$ton50.Range("E9").Select()

$ton420.Activate()
$ton420.Range("E5").Select()
